$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before C (shifts old Integral/Time columns to D/E)
$ws.Range("C:C").Insert()

# New column header: Variance
$ws.Range("C1").Value = "Variance"

# Fill column C with the squared values of column B (variance = STD^2)
$ws.Range("C2").Formula = "=B2^2"
$ws.Range("C3:C11").Formula = "=B3^2"

# Row 13: averages row - fill C13 with the average of the variance column
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"

# Row 14: STD row - remove the old B14/C14 STD formulas (they no longer apply
# now that B is raw data and C is Variance); keep D14/E14 (old C14/D14, shifted)
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# Row 15: new RMS row
$ws.Range("A15").Value = "RMS"
$ws.Range("B15").Formula = "=SQRT(C13)"

# Restore selection to match the committed workbook state
$ws.Range("B13").Select()
